$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Row 4: new script run entry (Login) ---
$ws.Range("A4").Formula = "'1"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Formula = "Login"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Formula = "'2024-04-30"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Formula = "EN"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "PASS"

# --- Row 11: clear old entry (Load or Pay credit card) ---
$ws.Range("A11:E11").ClearContents()

# --- Row 16: clear old entry (Check net worth statment chasrts) ---
$ws.Range("A16:E16").ClearContents()

# --- Row 17: Check wazin account details (FAIL) ---
$ws.Range("A17").Formula = "'14"
$ws.Range("B17").Formula = " Check wazin account details "
$ws.Range("C17").Formula = "'2024-04-30"
$ws.Range("D17").Formula = "EN"
$ws.Range("E17").Formula = " FAIL "

# --- Row 18: Between my account transfer (PASS) ---
$ws.Range("A18").Formula = "'15"
$ws.Range("B18").Formula = "Between my account transfer"
$ws.Range("C18").Formula = "'2024-04-30"
$ws.Range("D18").Formula = "EN"
$ws.Range("E18").Formula = "PASS"

# --- Row 19: Within riyad bank trnasfer (PASS) ---
$ws.Range("A19").Formula = "'16"
$ws.Range("B19").Formula = "Within riyad bank trnasfer"
$ws.Range("C19").Formula = "'2024-04-30"
$ws.Range("D19").Formula = "EN"
$ws.Range("E19").Formula = "PASS"

# --- Update selection to reflect last active cell ---
$ws.Range("E12").Select()
